$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value2 = 0.2307692307692308
$ws.Cells.Item(2, 3).Value2 = 0.493006993006993
$ws.Cells.Item(2, 10).Value2 = 0.01048951048951049
$ws.Cells.Item(2, 16).Value2 = 0.1888111888111888
$ws.Cells.Item(2, 19).Value2 = 0.07692307692307693
$ws.Cells.Item(3, 2).Value2 = 0.01379310344827586
$ws.Cells.Item(3, 3).Value2 = 0.01379310344827586
$ws.Cells.Item(3, 10).Value2 = 0.02068965517241379
$ws.Cells.Item(3, 16).Value2 = 0.7310344827586207
$ws.Cells.Item(3, 19).Value2 = 0.2206896551724138
$ws.Cells.Item(4, 10).Value2 = 0.1136363636363636
$ws.Cells.Item(4, 16).Value2 = 0.6590909090909091
$ws.Cells.Item(4, 19).Value2 = 0.2272727272727273
$ws.Cells.Item(6, 2).Value2 = 0.04149377593360996
$ws.Cells.Item(6, 4).Value2 = 0.008298755186721992
$ws.Cells.Item(6, 6).Value2 = 0.05809128630705394
$ws.Cells.Item(6, 10).Value2 = 0.2821576763485477
$ws.Cells.Item(6, 15).Value2 = 0.04149377593360996
$ws.Cells.Item(6, 17).Value2 = 0.1535269709543569
$ws.Cells.Item(6, 18).Value2 = 0.05809128630705394
$ws.Cells.Item(6, 19).Value2 = 0.3568464730290457
$ws.Cells.Item(7, 2).Value2 = 0.1158798283261803
$ws.Cells.Item(7, 4).Value2 = 0.03004291845493562
$ws.Cells.Item(7, 6).Value2 = 0.03862660944206009
$ws.Cells.Item(7, 10).Value2 = 0.1888412017167382
$ws.Cells.Item(7, 15).Value2 = 0.02575107296137339
$ws.Cells.Item(7, 17).Value2 = 0.1373390557939914
$ws.Cells.Item(7, 18).Value2 = 0.07296137339055794
$ws.Cells.Item(7, 19).Value2 = 0.3905579399141631
$ws.Cells.Item(8, 2).Value2 = 0.0871559633027523
$ws.Cells.Item(8, 4).Value2 = 0.01605504587155963
$ws.Cells.Item(8, 6).Value2 = 0.06651376146788991
$ws.Cells.Item(8, 10).Value2 = 0.09174311926605505
$ws.Cells.Item(8, 15).Value2 = 0.03211009174311927
$ws.Cells.Item(8, 17).Value2 = 0.1513761467889908
$ws.Cells.Item(8, 18).Value2 = 0.1009174311926606
$ws.Cells.Item(8, 19).Value2 = 0.4541284403669725
$ws.Cells.Item(9, 2).Value2 = 0.0481283422459893
$ws.Cells.Item(9, 4).Value2 = 0.0160427807486631
$ws.Cells.Item(9, 6).Value2 = 0.0374331550802139
$ws.Cells.Item(9, 10).Value2 = 0.1443850267379679
$ws.Cells.Item(9, 15).Value2 = 0.03208556149732621
$ws.Cells.Item(9, 17).Value2 = 0.1925133689839572
$ws.Cells.Item(9, 18).Value2 = 0.1336898395721925
$ws.Cells.Item(9, 19).Value2 = 0.3957219251336899
$ws.Cells.Item(10, 2).Value2 = 0.1038251366120219
$ws.Cells.Item(10, 4).Value2 = 0.02029664324746292
$ws.Cells.Item(10, 6).Value2 = 0.06713505074160812
$ws.Cells.Item(10, 10).Value2 = 0.1249024199843872
$ws.Cells.Item(10, 15).Value2 = 0.02185792349726776
$ws.Cells.Item(10, 17).Value2 = 0.1943793911007026
$ws.Cells.Item(10, 18).Value2 = 0.07103825136612021
$ws.Cells.Item(10, 19).Value2 = 0.3965651834504293
$ws.Cells.Item(11, 7).Value2 = 0.1782729805013928
$ws.Cells.Item(11, 10).Value2 = 0.0947075208913649
$ws.Cells.Item(11, 11).Value2 = 0.2033426183844011
$ws.Cells.Item(11, 12).Value2 = 0.5153203342618384
$ws.Cells.Item(11, 19).Value2 = 0.008356545961002786
$ws.Cells.Item(12, 7).Value2 = 0.7589743589743589
$ws.Cells.Item(12, 10).Value2 = 0.1846153846153846
$ws.Cells.Item(12, 12).Value2 = 0.005128205128205128
$ws.Cells.Item(12, 19).Value2 = 0.05128205128205128
$ws.Cells.Item(13, 7).Value2 = 0.6190476190476191
$ws.Cells.Item(13, 10).Value2 = 0.2619047619047619
$ws.Cells.Item(13, 19).Value2 = 0.119047619047619
$ws.Cells.Item(15, 6).Value2 = 0.03211009174311927
$ws.Cells.Item(15, 8).Value2 = 0.1559633027522936
$ws.Cells.Item(15, 9).Value2 = 0.05045871559633028
$ws.Cells.Item(15, 10).Value2 = 0.3394495412844037
$ws.Cells.Item(15, 11).Value2 = 0.01834862385321101
$ws.Cells.Item(15, 13).Value2 = 0.01376146788990826
$ws.Cells.Item(15, 15).Value2 = 0.06422018348623854
$ws.Cells.Item(15, 19).Value2 = 0.3256880733944954
$ws.Cells.Item(16, 6).Value2 = 0.05524861878453038
$ws.Cells.Item(16, 8).Value2 = 0.1546961325966851
$ws.Cells.Item(16, 9).Value2 = 0.04972375690607735
$ws.Cells.Item(16, 10).Value2 = 0.3977900552486188
$ws.Cells.Item(16, 11).Value2 = 0.1546961325966851
$ws.Cells.Item(16, 13).Value2 = 0.005524861878453038
$ws.Cells.Item(16, 15).Value2 = 0.04972375690607735
$ws.Cells.Item(16, 19).Value2 = 0.1325966850828729
$ws.Cells.Item(17, 6).Value2 = 0.03381642512077294
$ws.Cells.Item(17, 8).Value2 = 0.1714975845410628
$ws.Cells.Item(17, 9).Value2 = 0.1014492753623188
$ws.Cells.Item(17, 10).Value2 = 0.4033816425120773
$ws.Cells.Item(17, 11).Value2 = 0.09420289855072464
$ws.Cells.Item(17, 13).Value2 = 0.01932367149758454
$ws.Cells.Item(17, 14).Value2 = 0.002415458937198068
$ws.Cells.Item(17, 15).Value2 = 0.05072463768115942
$ws.Cells.Item(17, 19).Value2 = 0.1231884057971015
$ws.Cells.Item(18, 6).Value2 = 0.01025641025641026
$ws.Cells.Item(18, 8).Value2 = 0.1230769230769231
$ws.Cells.Item(18, 9).Value2 = 0.1076923076923077
$ws.Cells.Item(18, 10).Value2 = 0.4358974358974359
$ws.Cells.Item(18, 11).Value2 = 0.1128205128205128
$ws.Cells.Item(18, 13).Value2 = 0.02564102564102564
$ws.Cells.Item(18, 14).Value2 = 0.005128205128205128
$ws.Cells.Item(18, 15).Value2 = 0.04615384615384616
$ws.Cells.Item(18, 19).Value2 = 0.1333333333333333
$ws.Cells.Item(19, 6).Value2 = 0.02494497432134996
$ws.Cells.Item(19, 8).Value2 = 0.2010271460014673
$ws.Cells.Item(19, 9).Value2 = 0.07703595011005136
$ws.Cells.Item(19, 10).Value2 = 0.3514306676449009
$ws.Cells.Item(19, 11).Value2 = 0.1276595744680851
$ws.Cells.Item(19, 13).Value2 = 0.01907556859867938
$ws.Cells.Item(19, 14).Value2 = 0.0007336757153338225
$ws.Cells.Item(19, 15).Value2 = 0.06162876008804109
$ws.Cells.Item(19, 19).Value2 = 0.136463683052091
